$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.118.62"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.484.39"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.57"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.04"
$ws.Range("E6").Value = "  -5.56%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -2.94%  "
$ws.Range("D9").Value = "2.486.52"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("E13").Value = "  -4.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.95"
$ws.Range("E14").Value = "  -4.47%  "
$ws.Range("D15").Value = "2.936.61"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").Value = "66.916.47"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "2.517.41"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.74"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.80"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.34"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("B23").Value = "NEARProtocol"
$ws.Range("C23").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.42"
$ws.Range("E23").Value = "  -7.09%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.91"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  -6.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  -8.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "2.608.97"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("E30").Value = "  -6.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "507.46"
$ws.Range("E32").Value = "  -7.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.83"
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("E34").Value = "  -6.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.61"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.05"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.42"
$ws.Range("E39").Value = "  -3.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.56"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("E41").Value = "  -5.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.96"
$ws.Range("E42").Value = "  -5.84%  "
$ws.Range("E43").Value = "  -6.46%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -4.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.43"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.06"
$ws.Range("E47").Value = "  -4.47%  "
$ws.Range("E48").Value = "  -4.94%  "
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("E50").Value = "  -5.97%  "
$ws.Range("E51").Value = "  -4.33%  "
